$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update F18 cell content (was "Hoàn tất giao diện chỉnh sửa", now "Không có Rich text editor và FileUpload filter.")
$ws.Range("F18").Value = "Không có Rich text editor và FileUpload filter."

# Update E18 cell - was empty, now "X"
$ws.Range("E18").Value = "X"

# Row heights
$ws.Rows.Item(1).RowHeight = 47.25
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(18).RowHeight = 30

# View: scroll to show A7 at top-left, and select A17
$ws.Range("A17").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
